$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.423.73"
$ws.Range("E2").Value = "  -0.86%  "

$ws.Range("D3").Value = "1.851.15"
$ws.Range("E3").Value = "  +0.00%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6295"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.0000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07673"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.45%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2978"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.27%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.18%  "

$ws.Range("D11").Value = "1.971.97"
$ws.Range("E11").Value = "  +5.57%  "

$ws.Range("E12").Value = "  +1.14%  "

$ws.Range("E13").Value = "  -0.66%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6899"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.64%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001001"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.48%  "

$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.32%  "

$ws.Range("D17").Value = "2.193.28"
$ws.Range("E17").Value = "  +3.65%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.192"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.08%  "

$ws.Range("D19").Value = "29.546.23"
$ws.Range("E19").Value = "  -0.62%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "232.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.69%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9998"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.670"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.49%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9997"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.22%  "

$ws.Range("E25").Value = "  -2.10%  "

$ws.Range("E26").Value = "  -1.74%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.492"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.475"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05778"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.258"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.133"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.015"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.880"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7206"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.585"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.51%  "

$ws.Range("D38").Value = "1.251.49"
$ws.Range("E38").Value = "  +4.26%  "

$ws.Range("E39").Value = "  -0.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01810"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.37%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9080"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.091"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.06%  "

$ws.Range("D43").Value = "2.115.95"
$ws.Range("E43").Value = "  +4.56%  "

$ws.Range("E44").Value = "  -0.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "68.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.30%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.69"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.57%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.313"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000120"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.209"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4035"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.700"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.83%  "
